# Word COM-interop script implementing the diff:
#  - Replace the trailing empty paragraph (paraId 2A38E721) with three
#    paragraphs: the first keeps the paragraph mark and gets "9a ..." text
#    (losing its <w:u/> underline run property), and two new numbered
#    ("Paragrafoelenco" / numId 43) list paragraphs are appended after it.
#  - Add a brand-new numbering list definition (abstractNum "35", a
#    hybridMultilevel decimal list) ahead of the abstractNum that used to
#    be numbered "35", shifting every later abstractNum id up by one and
#    wiring a new <w:num numId="43"> to the freshly inserted abstractNum.

$d = $word.ActiveDocument

# Pull the whole package (all parts, including word/numbering.xml) as one
# WordprocessingML package string so document.xml and numbering.xml can be
# edited together and written back atomically via InsertXML.
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 1) numbering.xml : renumber abstractNum ids 35..41 -> 36..42
#    (processed through unique placeholders so the shifts can't collide
#    with each other, regardless of order)
# ---------------------------------------------------------------------
$shiftIds = 41, 40, 39, 38, 37, 36, 35
foreach ($oldId in $shiftIds) {
    $xml = $xml.Replace('<w:abstractNum w:abstractNumId="' + $oldId + '"', '<w:abstractNum w:abstractNumId="__SHIFT_' + $oldId + '__"')
}
foreach ($oldId in $shiftIds) {
    $newId = $oldId + 1
    $xml = $xml.Replace('__SHIFT_' + $oldId + '__', [string]$newId)
}

# ---------------------------------------------------------------------
# 2) numbering.xml : point the <w:num> entries that used to reference the
#    shifted abstractNum ids at their new ids (num 42 kept its historical
#    abstractNumId=35 definition, which now lives at abstractNumId=36).
# ---------------------------------------------------------------------
function Update-NumMapping($xmlText, $numId, $oldAbs, $newAbs) {
    $old = '<w:num w:numId="' + $numId + '"><w:abstractNumId w:val="' + $oldAbs + '"/></w:num>'
    $new = '<w:num w:numId="' + $numId + '"><w:abstractNumId w:val="' + $newAbs + '"/></w:num>'
    return $xmlText.Replace($old, $new)
}

$xml = Update-NumMapping $xml 15 41 42
$xml = Update-NumMapping $xml 24 40 41
$xml = Update-NumMapping $xml 28 36 37
$xml = Update-NumMapping $xml 31 38 39
$xml = Update-NumMapping $xml 33 37 38
$xml = Update-NumMapping $xml 34 39 40
$xml = Update-NumMapping $xml 42 35 36

# ---------------------------------------------------------------------
# 3) numbering.xml : insert the new abstractNum "35" definition right
#    before the (now renumbered) abstractNum "36", and append the new
#    <w:num numId="43"> entry that points at it, right before
#    </w:numbering>.
# ---------------------------------------------------------------------
$newAbstractNum = '<w:abstractNum w:abstractNumId="35" w15:restartNumberingAfterBreak="0"><w:nsid w:val="665045E8"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="49C6BDA8"/><w:lvl w:ilvl="0" w:tplc="8B00DFA4"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%1."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="705" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="04100019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%2."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1425" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="2" w:tplc="0410001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%3."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="2145" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="3" w:tplc="0410000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%4."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2865" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="4" w:tplc="04100019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%5."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3585" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="5" w:tplc="0410001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%6."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="4305" w:hanging="180"/></w:pPr></w:lvl><w:lvl w:ilvl="6" w:tplc="0410000F" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="decimal"/><w:lvlText w:val="%7."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5025" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="7" w:tplc="04100019" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerLetter"/><w:lvlText w:val="%8."/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5745" w:hanging="360"/></w:pPr></w:lvl><w:lvl w:ilvl="8" w:tplc="0410001B" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="lowerRoman"/><w:lvlText w:val="%9."/><w:lvlJc w:val="right"/><w:pPr><w:ind w:left="6465" w:hanging="180"/></w:pPr></w:lvl></w:abstractNum>'

$insertionMarker = '<w:abstractNum w:abstractNumId="36" w15:restartNumberingAfterBreak="0"><w:nsid w:val="668F45FE"/>'
$xml = $xml.Replace($insertionMarker, $newAbstractNum + $insertionMarker)

$newNumEntry = '<w:num w:numId="43"><w:abstractNumId w:val="35"/></w:num>'
$xml = $xml.Replace('</w:numbering>', $newNumEntry + '</w:numbering>')

# ---------------------------------------------------------------------
# 4) document.xml : replace the trailing empty paragraph with the new
#    text paragraph plus the two new numbered-list paragraphs.
# ---------------------------------------------------------------------
$oldParagraph = '<w:p w14:paraId="2A38E721" w14:textId="77777777" w:rsidR="00390F90" w:rsidRPr="0046674C" w:rsidRDefault="00390F90" w:rsidP="00FC41C4"><w:pPr><w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>'

$newParagraphs = '<w:p w14:paraId="2A38E721" w14:textId="77777777" w:rsidR="00390F90" w:rsidRPr="0046674C" w:rsidRDefault="00390F90" w:rsidP="00FC41C4"><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">9a </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> Il contadino vuole scegliere una specifica coltivazione</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="43"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>L' + [char]0x2019 + 'impiegato sceglie la coltivazione di interesse</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="43"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Il sistema mostra i dettagli relativi a quella coltivazione inclusi parametri ambientali target e attuali</w:t></w:r></w:p>'

$xml = $xml.Replace($oldParagraph, $newParagraphs)

# ---------------------------------------------------------------------
# 5) Write the modified package back into the document in one shot.
# ---------------------------------------------------------------------
$result = $d.Content.InsertXML($xml)

Write-Output "edit applied"
